$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The POD (column A) city names and POL (column B) terminal names are being
# re-stylized from ALL CAPS to Title Case. The underlying data values /
# row ordering are unchanged; only the text capitalization differs.
$cities = @("New York", "Savannah", "Miami", "Houston", "Indianapolis", "Los Angeles", "San Francisco")
$pols = @("Rotterdam", "Varna")

$row = 2
foreach ($pol in $pols) {
  foreach ($city in $cities) {
    for ($i = 0; $i -lt 5; $i++) {
      $ws.Cells.Item($row, 1).Value = $city
      $row = $row + 1
    }
  }
}

$row = 2
foreach ($pol in $pols) {
  foreach ($city in $cities) {
    for ($i = 0; $i -lt 5; $i++) {
      $ws.Cells.Item($row, 2).Value = $pol
      $row = $row + 1
    }
  }
}

# Update the active cell selection to match the post-edit workbook state.
$ws.Range("F11").Select()
